$d = $word.ActiveDocument

# The closing section of the letter template currently reads (right after
# the "=body" MERGEFIELD paragraph, up to the end of the document):
#   [empty]
#   "Yours sincerely"
#   [empty] x4
#   "Application Team"
#   "Data Protection Compliance Team"
#   "Ministry of Justice"
#   [empty]
#   [empty]
#
# Per the commit message ("Remove footer text from .docx letter template",
# with the 'from' footer text moving into the seeder instead) all of that
# fixed sign-off text is removed from the template, leaving a single
# trailing empty paragraph whose paragraph mark is now bold.

# Locate the paragraph that holds the "=body" MERGEFIELD result.
$bodyParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*body*") {
        $bodyParaIndex = $i
    }
}

# Everything from the paragraph right after "=body" through to the
# second-to-last paragraph of the document gets removed, leaving just the
# final (already-empty) trailing paragraph in place.
$startPara = $d.Paragraphs.Item($bodyParaIndex + 1)
$endPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)

$deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$deleteRange.Delete()

# The one remaining trailing paragraph's mark should become bold. A fully
# empty range can't carry character formatting on its own, so insert a
# placeholder character, bold the paragraph, then remove the character -
# the paragraph mark keeps the bold run property once the character is
# gone.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertAfter("X")
$lastPara.Range.Font.Bold = 1
$charRange = $d.Range($lastPara.Range.Start, $lastPara.Range.Start + 1)
$charRange.Delete()
